$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenario")

# Update SBW (column E) and Target Weight (column G) for the scenario rows 2-7
$ws.Range("E2:E7").Value = 365
$ws.Range("G2:G7").Value = 528

# Activate the sheet and move the selection/view back to E2:E7
$ws.Activate()
$ws.Range("E2:E7").Select()
